$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the last test-data row (row 6, the very long "abcdefghij..." FirstName row).
#    This shifts the used range from A1:E6 down to A1:E5.
$ws.Rows("6:6").Delete()

# 2. Rebuild the hyperlinks collection for the Email column (C2:C5) so the
#    stored relationship ids/targets line up with the new e-mail addresses
#    before we overwrite the cell text (Hyperlinks.Add also stamps display
#    text onto its anchor range, so this must happen before step 3).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Sarita2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3:C5"), "mailto:Sarita1@gmail.com", [Type]::Missing, [Type]::Missing, "Sarita1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Sarita2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Sarita3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Sarita4@gmail.com")

# 3. Update the Email (C) and mobile (D) columns with the new test data.
$ws.Range("C2").Value = "Sarita1@gmail.com"
$ws.Range("C3").Value = "Sarita2@gmail.com"
$ws.Range("C4").Value = "Sarita3@gmail.com"
$ws.Range("C5").Value = "Sarita4@gmail.com"

$ws.Range("D2").Value = "416-858-7771"
$ws.Range("D3").Value = "416-858-7772"
$ws.Range("D4").Value = "416-858-7773"
$ws.Range("D5").Value = "416-858-7774"

# 4. Resize the LastName/Email/mobile columns for the new (shorter) content.
$ws.Columns("B").ColumnWidth = 6.76
$ws.Columns("C").ColumnWidth = 20.42
$ws.Columns("D").ColumnWidth = 18.59

# 5. Move the active selection, matching the refreshed worksheet view.
$ws.Range("L6").Select() | Out-Null

# 6. Touch the page setup so the sheet carries explicit print/page settings.
$ws.PageSetup.Orientation = 1
